$wb = $excel.ActiveWorkbook

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Row 3: event was cancelled - update the title and mark the price as not for sale
    $ws.Range("C3").Value = "丽水·首届TCT国风动漫游戏嘉年华（取消）"
    $ws.Range("G3").Value = "不可售"

    # Row 5: interested-attendee count increased from 3 to 4
    $ws.Range("F5").Value = 4
}
